# Add a "BOD" (date-of-birth) column to the CustomerData report, between
# CIFName and AccType, filled with a constant value for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CustomerData")

# --- Insert the new column -------------------------------------------------
$ws.Columns("E:E").Insert()

# New column inherits column D's width/style (24.1640625, style index 3)
$ws.Columns("E:E").ColumnWidth = $ws.Columns("D:D").ColumnWidth

# Header + data for the new column
$ws.Range("E1").Value = "BOD"
$ws.Range("E2:E15").Value = "1991-09-03"

# --- Fix up the defined name that described the old (narrower) table ------
$nm = $wb.Names.Item("CustomerData!REPORT6_3")
$nm.RefersTo = "=CustomerData!`$B`$2:`$J`$11"

# --- Grow the remembered sort range to cover the new column ---------------
# (Sorted by the already-monotonic Seq column so the existing row order,
# which must stay untouched, is preserved exactly.)
$srt = $ws.Sort
$srt.SortFields.Clear()
$sf = $srt.SortFields.Add($ws.Range("A1:A11"))
Write-Output "sort field added"
$srt.SetRange($ws.Range("B1:K11"))
Write-Output "sort range set"
$srt.Apply()

# --- Selection matches the post-edit state ---------------------------------
$ws.Range("D11").Select()
